# Append 3 new daily rows (245-247) to the COVID-tracking sheet,
# mirroring the style/format of the last existing data row (244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 244

# Copy the formatting (style/number format) of the last data row's
# date cell onto the three new date cells so they keep the same
# "center/top aligned, bordered, date-formatted" style (s="2").
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A245:A247").PasteSpecial(-4122)

# Row 245: 2021-05-03
$ws.Cells.Item(245, 1).Value = 44319
$ws.Cells.Item(245, 2).Value = 0
$ws.Cells.Item(245, 3).Value = 0
$ws.Cells.Item(245, 4).Value = 0

# Row 246: 2021-05-04
$ws.Cells.Item(246, 1).Value = 44320
$ws.Cells.Item(246, 2).Value = 1
$ws.Cells.Item(246, 3).Value = 1
$ws.Cells.Item(246, 4).Value = 46.70714619336758

# Row 247: 2021-05-05
$ws.Cells.Item(247, 1).Value = 44321
$ws.Cells.Item(247, 2).Value = 0
$ws.Cells.Item(247, 3).Value = 1
$ws.Cells.Item(247, 4).Value = 46.70714619336758
